$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Visitor row (row 2)
$ws.Range("B2").Value = "Brookwood"
$ws.Range("C2").Value = "Broncos"
$ws.Range("D2").Value = "12-0"
$ws.Range("E2").Value = "OLD_GOLD"

# Home row (row 3)
$ws.Range("B3").Value = "Parkview"
$ws.Range("C3").Value = "Panthers"
$ws.Range("D3").Value = "0-12"
$ws.Range("E3").Value = "PURPLE"

# Update selection to match final state
$ws.Range("E2").Select()
